$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$dateTexts = @{
  2  = "01/01/2007"; 3  = "01/01/2007"; 4  = "01/01/2007"
  5  = "01/01/2008"; 6  = "01/01/2008"; 7  = "01/01/2008"
  8  = "01/01/2009"; 9  = "01/01/2009"; 10 = "01/01/2009"
  11 = "01/01/2010"; 12 = "01/01/2010"; 13 = "01/01/2010"
  14 = "01/01/2011"; 15 = "01/01/2011"; 16 = "01/01/2011"
  17 = "01/01/2012"; 18 = "01/01/2012"; 19 = "01/01/2012"
  20 = "01/01/2013"; 21 = "01/01/2013"; 22 = "01/01/2013"
  23 = "01/01/2014"; 24 = "01/01/2014"; 25 = "01/01/2014"
  26 = "01/01/2015"; 27 = "01/01/2015"; 28 = "01/01/2015"
  29 = "01/01/2016"; 30 = "01/01/2016"; 31 = "01/01/2016"
  32 = "01/01/2017"; 33 = "01/01/2017"; 34 = "01/01/2017"
  35 = "01/01/2018"; 36 = "01/01/2018"; 37 = "01/01/2018"
  38 = "01/01/2019"; 39 = "01/01/2019"; 40 = "01/01/2019"
  41 = "01/01/2020"; 42 = "01/01/2020"; 43 = "01/01/2020"
  44 = "01/01/2021"; 45 = "01/01/2021"; 46 = "01/01/2021"
  47 = "01/01/2022"; 48 = "01/01/2022"; 49 = "01/01/2022"
  50 = "01/01/2023"; 51 = "01/01/2023"; 52 = "01/01/2023"
}

foreach ($row in $dateTexts.Keys) {
  $cell = $ws.Cells.Item($row, 2)
  $cell.NumberFormat = "@"
  $cell.Value = $dateTexts[$row]
  $cell.Style = "Normal"
}
